# Add 2022-Q3 data
#
# 1) Summary sheet ("总计"): insert a new data row for "2022-Q3" at the top
#    of the data (row 2), pushing all existing quarters down by one row.
# 2) Insert a brand-new worksheet named "2022-Q3" right after "总计",
#    populated with the per-fund holdings for that quarter.

# Helper: write a text value into a cell. Excel's COM layer auto-detects
# numeric-looking strings (e.g. "37.71", "004475") and silently converts
# them to numbers on plain assignment, so a leading apostrophe is used to
# force a genuine text cell; the style is then reset to "Normal" so the
# quote-prefix flag doesn't leave a stray number-format behind.
function Set-TextCell($range, $value) {
    $range.Value2 = "'" + $value
    $range.Style = "Normal"
}

# Helper: write a numeric value into a cell.
function Set-NumCell($range, $value) {
    $range.Value2 = $value
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update the "总计" (summary) sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Give the new last row (row 9) the same formatting as the current last
# row (row 8, which carries the bold/border style used by column A)
# before we shift the data down.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

# Existing rows 2-8 shift down into rows 3-9 (read bottom-up so we never
# overwrite a value before it has been copied down). Column A is just a
# sequential 0-based row index, so it is set directly rather than copied.
for ($r = 8; $r -ge 2; $r--) {
    $target = $r + 1
    $seqIndex = $target - 2
    Set-NumCell $summary.Range("A$target") $seqIndex
    Set-TextCell $summary.Range("B$target") $summary.Range("B$r").Value2
    Set-NumCell $summary.Range("C$target") $summary.Range("C$r").Value2
    Set-NumCell $summary.Range("D$target") $summary.Range("D$r").Value2
}

# New row 2: 2022-Q3 totals.
Set-NumCell  $summary.Range("A2") 0
Set-TextCell $summary.Range("B2") "2022-Q3"
Set-NumCell  $summary.Range("C2") 8
Set-NumCell  $summary.Range("D2") 5.59

# ---------------------------------------------------------------------
# Step 2: insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
# Copy the "2022-Q2" sheet so the new sheet starts with identical layout,
# header row and formatting, then overwrite its data with the 2022-Q3
# figures.
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template has 12 fund rows (rows 2-13); 2022-Q3 only has 8 fund rows
# (rows 2-9), so drop the extra rows entirely.
$q3.Rows("10:13").Delete()

Set-NumCell  $q3.Range("A2") 0
Set-TextCell $q3.Range("B2") "004475"
Set-TextCell $q3.Range("C2") "华泰柏瑞富利灵活配置混合A"
Set-TextCell $q3.Range("D2") "37.71"
Set-TextCell $q3.Range("E2") "67.54"
Set-TextCell $q3.Range("F2") "6.74"
Set-TextCell $q3.Range("G2") "2.5417"
Set-NumCell  $q3.Range("H2") 2

Set-NumCell  $q3.Range("A3") 1
Set-TextCell $q3.Range("B3") "014597"
Set-TextCell $q3.Range("C3") "华泰柏瑞富利灵活配置混合C"
Set-TextCell $q3.Range("D3") "21.79"
Set-TextCell $q3.Range("E3") "67.54"
Set-TextCell $q3.Range("F3") "6.74"
Set-TextCell $q3.Range("G3") "1.4686"
Set-NumCell  $q3.Range("H3") 2

Set-NumCell  $q3.Range("A4") 2
Set-TextCell $q3.Range("B4") "002207"
Set-TextCell $q3.Range("C4") "前海开源金银珠宝主题精选混合C"
Set-TextCell $q3.Range("D4") "6.72"
Set-TextCell $q3.Range("E4") "90.85"
Set-TextCell $q3.Range("F4") "7.99"
Set-TextCell $q3.Range("G4") "0.5369"
Set-NumCell  $q3.Range("H4") 5

Set-NumCell  $q3.Range("A5") 3
Set-TextCell $q3.Range("B5") "001302"
Set-TextCell $q3.Range("C5") "前海开源金银珠宝主题精选混合A"
Set-TextCell $q3.Range("D5") "3.99"
Set-TextCell $q3.Range("E5") "90.85"
Set-TextCell $q3.Range("F5") "7.99"
Set-TextCell $q3.Range("G5") "0.3188"
Set-NumCell  $q3.Range("H5") 5

Set-NumCell  $q3.Range("A6") 4
Set-TextCell $q3.Range("B6") "003304"
Set-TextCell $q3.Range("C6") "前海开源沪港深核心资源灵活配置混合A"
Set-TextCell $q3.Range("D6") "3.45"
Set-TextCell $q3.Range("E6") "90.59"
Set-TextCell $q3.Range("F6") "8.01"
Set-TextCell $q3.Range("G6") "0.2763"
Set-NumCell  $q3.Range("H6") 3

Set-NumCell  $q3.Range("A7") 5
Set-TextCell $q3.Range("B7") "003175"
Set-TextCell $q3.Range("C7") "华泰柏瑞多策略灵活配置混合A"
Set-TextCell $q3.Range("D7") "3.32"
Set-TextCell $q3.Range("E7") "67.70"
Set-TextCell $q3.Range("F7") "6.67"
Set-TextCell $q3.Range("G7") "0.2214"
Set-NumCell  $q3.Range("H7") 2

Set-NumCell  $q3.Range("A8") 6
Set-TextCell $q3.Range("B8") "003305"
Set-TextCell $q3.Range("C8") "前海开源沪港深核心资源灵活配置混合C"
Set-TextCell $q3.Range("D8") "1.89"
Set-TextCell $q3.Range("E8") "90.59"
Set-TextCell $q3.Range("F8") "8.01"
Set-TextCell $q3.Range("G8") "0.1514"
Set-NumCell  $q3.Range("H8") 3

Set-NumCell  $q3.Range("A9") 7
Set-TextCell $q3.Range("B9") "015450"
Set-TextCell $q3.Range("C9") "华泰柏瑞多策略灵活配置混合C"
Set-TextCell $q3.Range("D9") "1.15"
Set-TextCell $q3.Range("E9") "67.70"
Set-TextCell $q3.Range("F9") "6.67"
Set-TextCell $q3.Range("G9") "0.0767"
Set-NumCell  $q3.Range("H9") 2
